# Atualização de bases das ligas, do dia: 20-06-2024 às 20:11
#
# For a handful of row pairs the match rows were swapped (row N now holds
# what used to be row N+1's data, and vice versa), while the leading "id"
# column (A) keeps its own original value. Swap columns B:AD between each
# pair of rows to reproduce the change.
#
# NOTE: nested arrays (an "array of pairs") get flattened by this host's
# PowerShell loop, so the row numbers are kept in two parallel flat arrays
# instead and walked with a classic for-loop.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsA = @(23, 29, 59, 128, 146, 149, 167, 218, 225, 238)
$rowsB = @(24, 30, 60, 129, 147, 150, 168, 219, 226, 239)

for ($i = 0; $i -lt $rowsA.Length; $i++) {
    $r1 = $rowsA[$i]
    $r2 = $rowsB[$i]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value = $vals2
    $range2.Value = $vals1
}
